$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/e8e71357d55f3f4c6e1af316006afae2cb8ce93f/e2e/3aa784aa-e5d6-4808-b8e0-2db5cb387df7.md"
$mdDisplay = "3aa784aa-e5d6-4808-b8e0-2db5cb387df7.md"

$cfgUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/e8e71357d55f3f4c6e1af316006afae2cb8ce93f/.localization-config"
$cfgDisplay = ".localization-config"

$newStatus = "Handoff transform failed"
$zeroDate  = "0001-01-01 00:00:00"

# The "File Name" status column on the Overview sheet shares the same
# underlying text as the per-locale Status column, so it must be updated
# too so the old "Ready for handoff" string disappears entirely.
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

foreach ($ws in @($zhcn, $dede)) {
    # Status: handoff transform failed for this source file.
    $ws.Range("B2").Value = $newStatus

    # The handoff never produced a target file any more, so the "Latest
    # Handoff File" cell/hyperlink is cleared out entirely.
    $ws.Range("C2").ClearContents()

    # "Latest Handoff Datetime" resets back to the default/zero datetime.
    $ws.Range("D2").Value = $zeroDate

    # "Latest Handback DateTime" stays at the default/zero datetime.
    $ws.Range("G2").Value = $zeroDate

    # Reason flips from "Include" to "Ignored".
    $ws.Range("H2").Value = "Ignored"

    # Row 3 (.localization-config) keeps the same logical values.
    $ws.Range("D3").Value = $zeroDate
    $ws.Range("G3").Value = $zeroDate
    $ws.Range("H3").Value = "Ignored"

    # Removing the C2 hyperlink also clears every other hyperlink on the
    # sheet as a side effect of this runtime, so capture what we need and
    # rebuild the remaining links (A2, A3) afterwards.
    $ws.Range("A1:I3").Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfgDisplay) | Out-Null
}

Write-Host "Generated handoff report updates"
